$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.317.88"
$ws.Range("E2").Value = "  -1.51%  "

$ws.Range("D3").Value = "1.905.23"
$ws.Range("E3").Value = "  -2.29%  "

$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "1.000"
$r.Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "334.21"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -2.13%  "

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "0.9999"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.4640"
$r.Style = "Normal"
$ws.Range("E7").Value = "  -2.76%  "

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.4144"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "48.00"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.08032"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -2.44%  "

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "1.018"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -1.74%  "

$ws.Range("E12").Value = "  -1.92%  "

$ws.Range("D13").Value = "1.941.46"
$ws.Range("E13").Value = "  -0.51%  "

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "5.970"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -2.91%  "

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "7.151"
$r.Style = "Normal"
$ws.Range("E15").Value = "  -3.13%  "

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "89.22"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -2.74%  "

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "1.000"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -0.22%  "

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "0.00001033"
$r.Style = "Normal"
$ws.Range("E18").Value = "  -2.38%  "

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "0.06580"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "17.68"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -1.94%  "

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "0.9986"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("D22").Value = "29.330.31"
$ws.Range("E22").Value = "  -1.36%  "

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "5.518"
$r.Style = "Normal"
$ws.Range("E23").Value = "  -1.29%  "

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "11.41"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +1.38%  "

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "2.198"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -3.39%  "

$ws.Range("D26").Value = "2.153.14"
$ws.Range("E26").Value = "  -1.25%  "

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "157.09"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -2.40%  "

$ws.Range("E28").Value = "  -1.73%  "

$ws.Range("E29").Value = "  -1.20%  "

$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "5.653"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -0.25%  "

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "117.21"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -4.63%  "

$ws.Range("E32").Value = "  +3.67%  "

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "0.09457"
$r.Style = "Normal"
$ws.Range("E33").Value = "  -1.88%  "

$ws.Range("E34").Value = "  -2.76%  "

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "3.551"
$r.Style = "Normal"
$ws.Range("E35").Value = "  -3.65%  "

$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "5.418"
$r.Style = "Normal"
$ws.Range("E36").Value = "  -1.52%  "

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.06119"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -2.36%  "

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.02258"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -2.45%  "

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "8.460"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "1.184"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -0.24%  "

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.5896"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -3.15%  "

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.9995"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -0.16%  "

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "10.22"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -4.76%  "

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.1834"
$r.Style = "Normal"
$ws.Range("E44").Value = "  -3.09%  "

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "2.377"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "1.262"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -1.37%  "

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.07505"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +1.48%  "

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.5571"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -2.37%  "

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "12.16"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -2.45%  "

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "1.931"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -2.77%  "

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "112.79"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +0.00%  "
